# Weekly refresh of the "Poroto granado" price series:
# a new week's record is inserted at row 10 (pushing the existing
# history down by one row, including the row that falls off the
# bottom of the table, which becomes the new last row, 133).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 10..132 down to 11..133, leaving a blank (but formatted) row 10.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with this week's record.
$ws.Cells.Item(10, 1).Value  = 8
$ws.Cells.Item(10, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(10, 3).Value  = "Coquimbo"
$ws.Cells.Item(10, 4).Value  = 45022
$ws.Cells.Item(10, 5).Value  = 4
$ws.Cells.Item(10, 6).Value  = 100112030
$ws.Cells.Item(10, 7).Value  = "Poroto granado"
$ws.Cells.Item(10, 8).Value  = "Sin especificar"
$ws.Cells.Item(10, 9).Value  = "Primera"
$ws.Cells.Item(10, 10).Value = 440
$ws.Cells.Item(10, 11).Value = 29000
$ws.Cells.Item(10, 12).Value = 30000
$ws.Cells.Item(10, 13).Value = 29500
$ws.Cells.Item(10, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 1180
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
